# Applies the "Last 2 paragraphs then spellcheck" edit:
#  - corrects a handful of yes/no typos in column E (the spellcheck pass)
#  - filters the sheet on column E = "no" (hides the "yes" rows)
#  - updates the view's selection and the stale _FilterDatabase defined name

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("features")

# --- Fix the yes/no ("spellcheck") values in column E ---------------------
$ws.Range("E6").Value2  = "no"
$ws.Range("E19").Value2 = "yes"
$ws.Range("E28").Value2 = "no"
$ws.Range("E29").Value2 = "no"
$ws.Range("E33").Value2 = "no"
$ws.Range("E47").Value2 = "no"
$ws.Range("E55").Value2 = "no"
$ws.Range("E56").Value2 = "no"
$ws.Range("E57").Value2 = "no"
$ws.Range("E61").Value2 = "no"
$ws.Range("E67").Value2 = "no"

# --- Apply an AutoFilter over the whole table, keeping only E = "no" ------
$rng = $ws.Range("A1:O71")
[void]$rng.AutoFilter(5, @("no"))

# --- Update the selection / view to match the final state -----------------
[void]$ws.Range("L71").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1

# --- Refresh the stale _xlnm._FilterDatabase defined name ------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=features!`$A`$1:`$O`$71"
    }
}

$wb.Save()
